$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad" / Changed date) holds the same serial date (45205,
# i.e. 2023-10-06) for every data row from row 2 through row 211. The
# workbook was refreshed and all these cells were bumped forward by one
# day to serial 45206 (2023-10-07).
for ($r = 2; $r -le 211; $r++) {
    $ws.Cells.Item($r, 3).Value = 45206
}
